$wb = $excel.ActiveWorkbook

$oldId = "b5db8074-2b82-4955-876a-2f2288b25e28"
$newId = "7d50da10-b19d-4fee-a522-92f93730b2b4"
$oldHash = "509171f16d2f7e2cafde1ea33cace57f9a19070d"
$newHash = "6f99343bcb121fa17461792acec233a80e6245f2"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
$wsOverview.Range("G2").Value = "2016-09-04 03:04:35"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newId.md"
$wsZh.Range("G2").Value = "$newId.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-04 03:04:30"
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("I2:J2").Style = "Normal"
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"
$wsZh.Columns("I:I").ColumnWidth = 18.6506053379604
$wsZh.Columns("J:J").ColumnWidth = 21.7054770333426

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newId.md"
$wsDe.Range("G2").Value = "$newId.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-09-04 03:04:35"
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("I2:J2").Style = "Normal"
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"
$wsDe.Columns("I:I").ColumnWidth = 18.6506053379604
$wsDe.Columns("J:J").ColumnWidth = 21.7054770333426
